$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.770179333333333
$ws.Range("H2").Value = 5.310538
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.324764666666667
$ws.Range("N2").Value = 6.974294
$ws.Range("O2").Value = 0.04473923998638302
$ws.Range("P2").Value = 0.04473923998638301
$ws.Range("Q2").Value = 4.115250367796889
$ws.Range("R2").Value = 37.03725331017201
$ws.Range("S2").Value = 0.04473923998638302
$ws.Range("T2").Value = 0.04473923998638301

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.770179333333333
$ws.Range("H3").Value = 5.310538
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.27491966666667
$ws.Range("N3").Value = 57.824759
$ws.Range("O3").Value = 0.3709387315842666
$ws.Range("P3").Value = 0.3709387315842665
$ws.Range("Q3").Value = 34.12006444559356
$ws.Range("R3").Value = 307.080580010342
$ws.Range("S3").Value = 0.3709387315842666
$ws.Range("T3").Value = 0.3709387315842665

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.770179333333333
$ws.Range("H4").Value = 5.310538
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.5843220284293504
$ws.Range("P4").Value = 0.5843220284293504
$ws.Range("Q4").Value = 53.74770432259444
$ws.Range("R4").Value = 483.72933890335
$ws.Range("S4").Value = 0.5843220284293504
$ws.Range("T4").Value = 0.5843220284293504
